# Bridge table (Model) changed: USERHOME bridge table now uses
# email / Home_name (varchar columns) instead of UserPk / HomePk (INT).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: HomePk -> Home_name (new shared string added first)
$ws.Range("A39").Value = "Home_name"

# Row 38: INT -> varchar (60)
$ws.Range("B38").Value = "varchar (60)"

# Row 39: INT -> varchar (20)
$ws.Range("B39").Value = "varchar (20)"

# Row 38: UserPk -> email (reuses existing shared string)
$ws.Range("A38").Value = "email"

# Update the view state to match where the author ended up editing.
$ws.Range("B40").Select()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 30
    $win.ScrollColumn = 1
}
